$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text in cell E8 ("Good Morning" -> "GIT UPDATE")
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active selection being on E8 (as seen in the saved sheetView)
$ws.Range("E8").Select()
